$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row style reference not needed - header unchanged

$ws.Cells.Item(2, 1).Value = 'Rohan Chowla'
$ws.Cells.Item(2, 2).Value = 14
$ws.Cells.Item(2, 3).Value = 5
$ws.Cells.Item(2, 4).Value = 6
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 34

$ws.Cells.Item(3, 1).Value = 'Kevin Lee'
$ws.Cells.Item(3, 2).Value = 12
$ws.Cells.Item(3, 3).Value = 4
$ws.Cells.Item(3, 4).Value = 5
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 28

$ws.Cells.Item(4, 1).Value = 'Roman Ramirez'
$ws.Cells.Item(4, 2).Value = 6
$ws.Cells.Item(4, 3).Value = 9
$ws.Cells.Item(4, 4).Value = 5
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 22

$ws.Cells.Item(5, 1).Value = 'Jason Jackson'
$ws.Cells.Item(5, 2).Value = 5
$ws.Cells.Item(5, 3).Value = 3
$ws.Cells.Item(5, 4).Value = 7
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 15

$ws.Cells.Item(6, 1).Value = 'Kevin Cooper'
$ws.Cells.Item(6, 2).Value = 4
$ws.Cells.Item(6, 3).Value = 6
$ws.Cells.Item(6, 4).Value = 5
$ws.Cells.Item(6, 5).Value = 0
$ws.Cells.Item(6, 6).Value = 14

$ws.Cells.Item(7, 1).Value = 'Aaron Carter'
$ws.Cells.Item(7, 2).Value = 3
$ws.Cells.Item(7, 3).Value = 8
$ws.Cells.Item(7, 4).Value = 7
$ws.Cells.Item(7, 5).Value = 0
$ws.Cells.Item(7, 6).Value = 14

$ws.Cells.Item(8, 1).Value = 'Coby Lovelace'
$ws.Cells.Item(8, 2).Value = 4
$ws.Cells.Item(8, 3).Value = 3
$ws.Cells.Item(8, 4).Value = 5
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 13

$ws.Cells.Item(9, 1).Value = 'Nathan Snow'
$ws.Cells.Item(9, 2).Value = 4
$ws.Cells.Item(9, 3).Value = 2
$ws.Cells.Item(9, 4).Value = 3
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 11

$ws.Cells.Item(10, 1).Value = 'Cason Duszak'
$ws.Cells.Item(10, 2).Value = 4
$ws.Cells.Item(10, 3).Value = 2
$ws.Cells.Item(10, 4).Value = 5
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 11

$ws.Cells.Item(11, 1).Value = 'Leah Baetcke'
$ws.Cells.Item(11, 2).Value = 4
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(11, 4).Value = 6
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 10

$ws.Cells.Item(12, 1).Value = 'Gabe Silverstein'
$ws.Cells.Item(12, 2).Value = 2
$ws.Cells.Item(12, 3).Value = 6
$ws.Cells.Item(12, 4).Value = 8
$ws.Cells.Item(12, 5).Value = 0
$ws.Cells.Item(12, 6).Value = 10

$ws.Cells.Item(13, 1).Value = 'Will Simpson'
$ws.Cells.Item(13, 2).Value = 2
$ws.Cells.Item(13, 3).Value = 5
$ws.Cells.Item(13, 4).Value = 5
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = 10

$ws.Cells.Item(14, 1).Value = 'Jack Massingill'
$ws.Cells.Item(14, 2).Value = 2
$ws.Cells.Item(14, 3).Value = 5
$ws.Cells.Item(14, 4).Value = 8
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(14, 6).Value = 9

$ws.Cells.Item(15, 1).Value = 'Eric LastName'
$ws.Cells.Item(15, 2).Value = 2
$ws.Cells.Item(15, 3).Value = 1
$ws.Cells.Item(15, 4).Value = 1
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 6

$ws.Cells.Item(16, 1).Value = 'Matthew Rusten'
$ws.Cells.Item(16, 2).Value = 2
$ws.Cells.Item(16, 3).Value = 1
$ws.Cells.Item(16, 4).Value = 3
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 6

$ws.Cells.Item(17, 1).Value = 'Carla Betancourt'
$ws.Cells.Item(17, 2).Value = 2
$ws.Cells.Item(17, 3).Value = 1
$ws.Cells.Item(17, 4).Value = 2
$ws.Cells.Item(17, 5).Value = 0
$ws.Cells.Item(17, 6).Value = 5

$ws.Cells.Item(18, 1).Value = 'Yvonne Nguyen'
$ws.Cells.Item(18, 2).Value = 1
$ws.Cells.Item(18, 3).Value = 3
$ws.Cells.Item(18, 4).Value = 4
$ws.Cells.Item(18, 5).Value = 0
$ws.Cells.Item(18, 6).Value = 5

$ws.Cells.Item(19, 1).Value = 'Ann Hall'
$ws.Cells.Item(19, 2).Value = 0
$ws.Cells.Item(19, 3).Value = 4
$ws.Cells.Item(19, 4).Value = 5
$ws.Cells.Item(19, 5).Value = 1
$ws.Cells.Item(19, 6).Value = 5

$ws.Cells.Item(20, 1).Value = 'Rose Roché'
$ws.Cells.Item(20, 2).Value = 1
$ws.Cells.Item(20, 3).Value = 2
$ws.Cells.Item(20, 4).Value = 6
$ws.Cells.Item(20, 5).Value = 0
$ws.Cells.Item(20, 6).Value = 4

$ws.Cells.Item(21, 1).Value = 'Luci Nguyen'
$ws.Cells.Item(21, 2).Value = 1
$ws.Cells.Item(21, 3).Value = 1
$ws.Cells.Item(21, 4).Value = 1
$ws.Cells.Item(21, 5).Value = 1
$ws.Cells.Item(21, 6).Value = 4

$ws.Cells.Item(22, 1).Value = 'Helen Dunn'
$ws.Cells.Item(22, 2).Value = 1
$ws.Cells.Item(22, 3).Value = 1
$ws.Cells.Item(22, 4).Value = 1
$ws.Cells.Item(22, 5).Value = 1
$ws.Cells.Item(22, 6).Value = 3

$ws.Cells.Item(23, 1).Value = 'Noah Dale'
$ws.Cells.Item(23, 2).Value = 1
$ws.Cells.Item(23, 3).Value = 0
$ws.Cells.Item(23, 4).Value = 3
$ws.Cells.Item(23, 5).Value = 1
$ws.Cells.Item(23, 6).Value = 3

$ws.Cells.Item(24, 1).Value = 'Kristian Banlaoi'
$ws.Cells.Item(24, 2).Value = 1
$ws.Cells.Item(24, 3).Value = 0
$ws.Cells.Item(24, 4).Value = 2
$ws.Cells.Item(24, 5).Value = 0
$ws.Cells.Item(24, 6).Value = 2

$ws.Cells.Item(25, 1).Value = 'Piper Parker'
$ws.Cells.Item(25, 2).Value = 1
$ws.Cells.Item(25, 3).Value = 0
$ws.Cells.Item(25, 4).Value = 4
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 6).Value = 2

$ws.Cells.Item(26, 1).Value = 'Reagan Fryatt'
$ws.Cells.Item(26, 2).Value = 0
$ws.Cells.Item(26, 3).Value = 2
$ws.Cells.Item(26, 4).Value = 2
$ws.Cells.Item(26, 5).Value = 0
$ws.Cells.Item(26, 6).Value = 2

$ws.Cells.Item(27, 1).Value = 'Anna Brown'
$ws.Cells.Item(27, 2).Value = 0
$ws.Cells.Item(27, 3).Value = 1
$ws.Cells.Item(27, 4).Value = 2
$ws.Cells.Item(27, 5).Value = 0
$ws.Cells.Item(27, 6).Value = 1

$ws.Cells.Item(28, 1).Value = 'Alex LastName'
$ws.Cells.Item(28, 2).Value = 0
$ws.Cells.Item(28, 3).Value = 1
$ws.Cells.Item(28, 4).Value = 2
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 6).Value = 1

$ws.Cells.Item(29, 1).Value = 'Julie Jackson'
$ws.Cells.Item(29, 2).Value = 0
$ws.Cells.Item(29, 3).Value = 1
$ws.Cells.Item(29, 4).Value = 2
$ws.Cells.Item(29, 5).Value = 0
$ws.Cells.Item(29, 6).Value = 1

$ws.Cells.Item(30, 1).Value = 'Carolyn LastName'
$ws.Cells.Item(30, 2).Value = 0
$ws.Cells.Item(30, 3).Value = 1
$ws.Cells.Item(30, 4).Value = 2
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(30, 6).Value = 1

$ws.Cells.Item(31, 1).Value = 'Brian Tafazoli'
$ws.Cells.Item(31, 2).Value = 0
$ws.Cells.Item(31, 3).Value = 0
$ws.Cells.Item(31, 4).Value = 2
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(31, 6).Value = 0

$ws.Cells.Item(32, 1).Value = 'Sam Tellis'
$ws.Cells.Item(32, 2).Value = 0
$ws.Cells.Item(32, 3).Value = 0
$ws.Cells.Item(32, 4).Value = 2
$ws.Cells.Item(32, 5).Value = 0
$ws.Cells.Item(32, 6).Value = 0

$ws.Cells.Item(33, 1).Value = 'Cassie Deering'
$ws.Cells.Item(33, 2).Value = 0
$ws.Cells.Item(33, 3).Value = 0
$ws.Cells.Item(33, 4).Value = 2
$ws.Cells.Item(33, 5).Value = 0
$ws.Cells.Item(33, 6).Value = 0

$ws.Cells.Item(34, 1).Value = 'Yafu LastName'
$ws.Cells.Item(34, 2).Value = 0
$ws.Cells.Item(34, 3).Value = 0
$ws.Cells.Item(34, 4).Value = 2
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(34, 6).Value = 0

$ws.Cells.Item(35, 1).Value = 'Kim LastName'
$ws.Cells.Item(35, 2).Value = 0
$ws.Cells.Item(35, 3).Value = 0
$ws.Cells.Item(35, 4).Value = 2
$ws.Cells.Item(35, 5).Value = 0
$ws.Cells.Item(35, 6).Value = 0

$ws.Cells.Item(36, 1).Value = 'Evan Sooklal'
$ws.Cells.Item(36, 2).Value = 0
$ws.Cells.Item(36, 3).Value = 0
$ws.Cells.Item(36, 4).Value = 4
$ws.Cells.Item(36, 5).Value = 0
$ws.Cells.Item(36, 6).Value = 0

$ws.Cells.Item(37, 1).Value = 'Paul Bartenfeld'
$ws.Cells.Item(37, 2).Value = 0
$ws.Cells.Item(37, 3).Value = 0
$ws.Cells.Item(37, 4).Value = 8
$ws.Cells.Item(37, 5).Value = 0
$ws.Cells.Item(37, 6).Value = 0

# Apply the same style to the newly added rows (30-37) column A as the rest of column A
$srcStyleRange = $ws.Range("A29")
$srcStyleRange.Copy()
$dstStyleRange = $ws.Range("A30:A37")
$dstStyleRange.PasteSpecial(-4122)
$excel.CutCopyMode = 0
